$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.217.79'
$ws.Range("E2").Value = '  -5.15%  '
$ws.Range("D3").Value = '3.229.83'
$ws.Range("E3").Value = '  -8.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.26'
$ws.Range("E5").Value = '  -5.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.62'
$ws.Range("E6").Value = '  -12.11%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.224.71'
$ws.Range("E8").Value = '  -8.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.543'
$ws.Range("E9").Value = '  -11.17%  '
$ws.Range("E10").Value = '  -13.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.70'
$ws.Range("E11").Value = '  -7.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.503'
$ws.Range("E12").Value = '  -14.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.26'
$ws.Range("E13").Value = '  -18.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000244'
$ws.Range("E14").Value = '  -11.86%  '
$ws.Range("D15").Value = '3.749.51'
$ws.Range("E15").Value = '  -8.66%  '
$ws.Range("D16").Value = '67.122.91'
$ws.Range("E16").Value = '  -5.32%  '
$ws.Range("D17").Value = '3.228.31'
$ws.Range("E17").Value = '  -8.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '541.58'
$ws.Range("E18").Value = '  -11.60%  '
$ws.Range("E19").Value = '  -5.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.15'
$ws.Range("E20").Value = '  -15.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.12'
$ws.Range("E21").Value = '  -15.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.759'
$ws.Range("E22").Value = '  -14.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.78'
$ws.Range("E23").Value = '  -13.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.59'
$ws.Range("E24").Value = '  -12.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.48'
$ws.Range("E25").Value = '  -14.53%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  -16.65%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.07'
$ws.Range("E28").Value = '  -12.09%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '29.41'
$ws.Range("E29").Value = '  -13.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.14'
$ws.Range("E30").Value = '  -17.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.66'
$ws.Range("E31").Value = '  -12.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.14'
$ws.Range("E32").Value = '  -13.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '544.14'
$ws.Range("E33").Value = '  -10.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.54'
$ws.Range("E34").Value = '  -20.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.71'
$ws.Range("E35").Value = '  -16.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.51'
$ws.Range("E37").Value = '  -6.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0440'
$ws.Range("E38").Value = '  -7.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.18'
$ws.Range("E39").Value = '  -15.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0845'
$ws.Range("E40").Value = '  -16.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.127'
$ws.Range("E41").Value = '  -13.15%  '
$ws.Range("D42").Value = '2.923.43'
$ws.Range("E42").Value = '  -13.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.57'
$ws.Range("E43").Value = '  -27.18%  '
$ws.Range("E44").Value = '  -21.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.260'
$ws.Range("E45").Value = '  -17.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.38'
$ws.Range("E46").Value = '  -20.01%  '
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.06'
$ws.Range("E48").Value = '  -19.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.11'
$ws.Range("E49").Value = '  -18.25%  '
$ws.Range("E50").Value = '  -13.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '123.57'
$ws.Range("E51").Value = '  -7.48%  '
